$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B5 value
$ws.Range("B5").Value = 80221

# Remove C5 (was "Ovaliderad")
$ws.Range("C5").ClearContents()

# Update Q5 and R5 values
$ws.Range("Q5").Value = 551093
$ws.Range("R5").Value = 7027429

# Remove Z5 and AB5 (were "00:00")
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
